# Add the new monthly data row (2024-08-01 / serial 45505) to the bottom of
# the inflation-contribution-weights table, and switch the date column's
# display format from "yyyy-mm-dd hh:mm:ss" to "yyyy/mm/dd hh:mm:ss".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Append the new row of data (row 20) ---
$ws.Range("A20").Value = 45505
$ws.Range("B20").Value = 0.24178
$ws.Range("C20").Value = 0.21784
$ws.Range("D20").Value = 0.23053
$ws.Range("E20").Value = 0.30989
$ws.Range("F20").Value = 0.25091

# Give the new date cell the same formatting (font/number format) as the
# rest of the date column before restyling, so the whole column consolidates
# onto a single shared style.
$ws.Range("A19").Copy()
$ws.Range("A20").PasteSpecial(-4122)

# --- 2. Change the date format for the whole date column (A2:A20) ---
$ws.Range("A2:A20").NumberFormat = "yyyy/mm/dd hh:mm:ss"
